$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the COUNTIF range bound from A1:A10000 to A1:A10005 for each incident-type
# count formula (loitering, littering, noise complaint, land argument,
# discrimination, domestic abuse, assault, multiple assaults, traffic incident,
# fatality).
$ws.Range("B2").Formula = '=COUNTIF(A1:A10005,"*LOITERING_TEXT_*")'
$ws.Range("B6").Formula = '=COUNTIF(A1:A10005,"*LITTERING_TEXT_*")'
$ws.Range("B11").Formula = '=COUNTIF(A1:A10005,"*NOISECOMPLAINT_TEXT_*")'
$ws.Range("B18").Formula = '=COUNTIF(A1:A10005,"*LANDARGUMENT_TEXT_*")'
$ws.Range("B22").Formula = '=COUNTIF(A1:A10005,"*DISCRIMINATION_TEXT_*")'
$ws.Range("B24").Formula = '=COUNTIF(A1:A10005,"*DOMESTICABUSE_TEXT_*")'
$ws.Range("B27").Formula = '=COUNTIF(A1:A10005,"*ASSAULT_TEXT_*")'
$ws.Range("B33").Formula = '=COUNTIF(A1:A10005,"*MULTIPLEASSAULTS_TEXT_*")'
$ws.Range("B37").Formula = '=COUNTIF(A1:A10005,"*TRAFFICINCIDENT_TEXT_*")'
$ws.Range("B42").Formula = '=COUNTIF(A1:A10005,"*FATALITY_TEXT_*")'

# Update the localization string table: fix the line-break marker in
# START_SCREEN_DESCRIPTION's body text (\n -> *n*) and append all of the new
# localization keys/values used to hook up the rest of the game UI.
$ws.Range("A45").Value = 'START_SCREEN_DESCRIPTION'
$ws.Range("B45").Value = 'Welcome to resource force, take on the role of a police operator and help to solve crimes across the city. You must use your resources carefully to solve cases and keep civilians happy and feeling safe.*n*But be careful, officers are limited so use them wisely'
$ws.Range("A46").Value = 'APP_NAME'
$ws.Range("B46").Value = 'ResourceForce'
$ws.Range("A47").Value = 'START_SCREEN_TAP'
$ws.Range("B47").Value = 'Tap to start!'
$ws.Range("A48").Value = 'BASIC_TEXT_OK'
$ws.Range("B48").Value = 'OK!'
$ws.Range("C48").Value = 'XXXX'
$ws.Range("D48").Value = 'XXXX'
$ws.Range("E48").Value = 'XXXX'
$ws.Range("A49").Value = 'BASIC_TEXT_WAIT'
$ws.Range("B49").Value = 'Wait!'
$ws.Range("C49").Value = 'XXXX'
$ws.Range("D49").Value = 'XXXX'
$ws.Range("E49").Value = 'XXXX'
$ws.Range("A50").Value = 'BASIC_TEXT_SEND_ONE'
$ws.Range("B50").Value = 'Send 1 officer for {0} turns'
$ws.Range("C50").Value = 'XXXX'
$ws.Range("D50").Value = 'XXXX'
$ws.Range("E50").Value = 'XXXX'
$ws.Range("A51").Value = 'BASIC_TEXT_SEND_MANY'
$ws.Range("B51").Value = 'Send {0} officers for {1} turns'
$ws.Range("C51").Value = 'XXXX'
$ws.Range("D51").Value = 'XXXX'
$ws.Range("E51").Value = 'XXXX'
$ws.Range("A52").Value = 'BASIC_TEXT_DEVELOPED'
$ws.Range("B52").Value = 'Developed Case'
$ws.Range("C52").Value = 'XXXX'
$ws.Range("D52").Value = 'XXXX'
$ws.Range("E52").Value = 'XXXX'
$ws.Range("A53").Value = 'BASIC_TEXT_ARREST_SUCCESS'
$ws.Range("B53").Value = 'Arrests have been made'
$ws.Range("C53").Value = 'XXXX'
$ws.Range("D53").Value = 'XXXX'
$ws.Range("E53").Value = 'XXXX'
$ws.Range("A54").Value = 'BASIC_TEXT_ARREST_FAIL'
$ws.Range("B54").Value = 'Officers fail to make any arrests regarding the case'
$ws.Range("C54").Value = 'XXXX'
$ws.Range("D54").Value = 'XXXX'
$ws.Range("E54").Value = 'XXXX'
$ws.Range("A55").Value = 'BASIC_TEXT_ASK_CITIZEN'
$ws.Range("B55").Value = 'Ask citizen for help'
$ws.Range("C55").Value = 'XXXX'
$ws.Range("D55").Value = 'XXXX'
$ws.Range("E55").Value = 'XXXX'
$ws.Range("A56").Value = 'BASIC_TEXT_CITIZEN_SUCCESS'
$ws.Range("B56").Value = 'Citizens provide evidence through the INSPEC2T app, 2 have been charged'
$ws.Range("C56").Value = 'XXXX'
$ws.Range("D56").Value = 'XXXX'
$ws.Range("E56").Value = 'XXXX'
$ws.Range("A57").Value = 'BASIC_TEXT_CITIZEN_FAIL'
$ws.Range("B57").Value = 'Citizen fails to provide any solid evidence for the case'
$ws.Range("C57").Value = 'XXXX'
$ws.Range("D57").Value = 'XXXX'
$ws.Range("E57").Value = 'XXXX'
$ws.Range("A58").Value = 'BASIC_TEXT_NO_MORE_INCIDENTS'
$ws.Range("B58").Value = 'No more incidents to check this turn'
$ws.Range("C58").Value = 'XXXX'
$ws.Range("D58").Value = 'XXXX'
$ws.Range("E58").Value = 'XXXX'
$ws.Range("A59").Value = 'BASIC_TEXT_NEXT_TURN'
$ws.Range("B59").Value = 'Next Turn'
$ws.Range("C59").Value = 'XXXX'
$ws.Range("D59").Value = 'XXXX'
$ws.Range("E59").Value = 'XXXX'
$ws.Range("A60").Value = 'INCIDENT_NEW'
$ws.Range("B60").Value = 'New'
$ws.Range("C60").Value = 'XXXX'
$ws.Range("D60").Value = 'XXXX'
$ws.Range("E60").Value = 'XXXX'
$ws.Range("A61").Value = 'INCIDENT_ONGOING'
$ws.Range("B61").Value = 'Ongoing'
$ws.Range("C61").Value = 'XXXX'
$ws.Range("D61").Value = 'XXXX'
$ws.Range("E61").Value = 'XXXX'
$ws.Range("A62").Value = 'INCIDENT_RESOLVED'
$ws.Range("B62").Value = 'Resolved'
$ws.Range("C62").Value = 'XXXX'
$ws.Range("D62").Value = 'XXXX'
$ws.Range("E62").Value = 'XXXX'
$ws.Range("A63").Value = 'INCIDENT_CASE'
$ws.Range("B63").Value = 'Case'
$ws.Range("C63").Value = 'XXXX'
$ws.Range("D63").Value = 'XXXX'
$ws.Range("E63").Value = 'XXXX'
$ws.Range("A64").Value = 'INCIDENT_CASE_SUBJECT'
$ws.Range("B64").Value = 'Subject'
$ws.Range("C64").Value = 'XXXX'
$ws.Range("D64").Value = 'XXXX'
$ws.Range("E64").Value = 'XXXX'
$ws.Range("A65").Value = 'INCIDENT_OFFICERS'
$ws.Range("B65").Value = 'Officers'
$ws.Range("C65").Value = 'XXXX'
$ws.Range("D65").Value = 'XXXX'
$ws.Range("E65").Value = 'XXXX'
$ws.Range("A66").Value = 'BASIC_TEXT_RESOLVED_CASES'
$ws.Range("B66").Value = 'Arrests'
$ws.Range("C66").Value = 'XXXX'
$ws.Range("D66").Value = 'XXXX'
$ws.Range("E66").Value = 'XXXX'
$ws.Range("A67").Value = 'BASIC_TEXT_ACTIVE_CASES'
$ws.Range("B67").Value = 'Active Cases'
$ws.Range("C67").Value = 'XXXX'
$ws.Range("D67").Value = 'XXXX'
$ws.Range("E67").Value = 'XXXX'
$ws.Range("A68").Value = 'BASIC_TEXT_TURN'
$ws.Range("B68").Value = 'Turn'
$ws.Range("C68").Value = 'XXXX'
$ws.Range("D68").Value = 'XXXX'
$ws.Range("E68").Value = 'XXXX'
$ws.Range("A69").Value = 'BASIC_TEXT_GAMEOVER'
$ws.Range("B69").Value = 'Game Over*n*Too many unresolved cases'
$ws.Range("C69").Value = 'XXXX'
$ws.Range("D69").Value = 'XXXX'
$ws.Range("E69").Value = 'XXXX'
$ws.Range("A70").Value = 'BASIC_TEXT_GAMEOVER_BODY'
$ws.Range("B70").Value = 'You Survived {0} Turns*n*And Made Arrests for {1}% of Cases'
$ws.Range("C70").Value = 'XXXX'
$ws.Range("D70").Value = 'XXXX'
$ws.Range("E70").Value = 'XXXX'
$ws.Range("A71").Value = 'BASIC_TEXT_AVAILABLE'
$ws.Range("B71").Value = 'Available'
$ws.Range("C71").Value = 'XXXX'
$ws.Range("D71").Value = 'XXXX'
$ws.Range("E71").Value = 'XXXX'
$ws.Range("A72").Value = 'BASIC_TEXT_TURNS_UNTIL_AVAILABLE'
$ws.Range("B72").Value = 'turns until available'
$ws.Range("C72").Value = 'XXXX'
$ws.Range("D72").Value = 'XXXX'
$ws.Range("E72").Value = 'XXXX'

# Reflect the new extent of the table and the view/selection state that was
# active when the additional rows were appended.
$win = $excel.ActiveWindow
$win.ScrollRow = 34
$win.ScrollColumn = 1
$ws.Range("B70").Select()

